$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before existing row 2, shifting old data down
$ws.Rows("2:3").Insert()

# New row 2
$ws.Range("A2").Value = -83.31174
$ws.Range("B2").Value = 9.0062099999999994

# New row 3
$ws.Range("A3").Value = -83.444000000000003
$ws.Range("B3").Value = 9.8290000000000006

# Update selection to D5 (as in the diff)
$ws.Range("D5").Select()
